$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the title/header block strings in column E (rows 2-8).
$ws.Range("E2:E8").ClearContents()

# Clear the helper "category" labels in column H (rows 11-54) that are no
# longer used now that the report is generated directly off column C.
$ws.Range("H11:H23").ClearContents()
$ws.Range("H26:H54").ClearContents()

# Clear the balancing-check formula that depended on the helper column.
$ws.Range("H56").ClearContents()

# Restore the active selection to match the regenerated report state.
$ws.Range("H1:H1048576").Select()
